$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the style that the current last row (row 47, col A) uses,
# so the new last row (row 48) can reuse it.
$lastRowDateFormat = $ws.Range("A47").NumberFormat

# Row 47 is no longer the last row - give it the same format as the
# rest of the data rows (copy from A46, a regular interior row).
$ws.Range("A47").NumberFormat = $ws.Range("A46").NumberFormat

# Append the new daily row (48) with the next day's data.
$ws.Range("A48").Value = 45788
$ws.Range("B48").Value = 196
$ws.Range("C48").Value = 207
$ws.Range("D48").Value = 204

# The new last row gets the "last row" date-only formatting that row 47
# used to have.
$ws.Range("A48").NumberFormat = $lastRowDateFormat
